$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 99
for ($r = 2; $r -le $lastRow; $r++) {
    $dCell = $ws.Cells.Item($r, 4)
    $eCell = $ws.Cells.Item($r, 5)
    $fCell = $ws.Cells.Item($r, 6)

    $d = $dCell.Value2
    $e = $eCell.Value2
    $f = $fCell.Value2

    if ($d -eq $null -or $e -eq $null -or $f -eq $null) {
        continue
    }

    try {
        $fDate = [datetime]::ParseExact([string]$f, "yyyyMMdd", $null)
    } catch {
        # Unparseable start date (e.g. malformed value) - leave this row untouched.
        continue
    }

    if ($e -le 1) {
        # Cycle complete: reset remaining days to the full cycle length and
        # roll the start date forward by that many days.
        $newE = $d
        $newDate = $fDate.AddDays([double]$d)
        $newF = [int]$newDate.ToString("yyyyMMdd")

        $eCell.Value2 = $newE
        $fCell.Value2 = $newF
    } else {
        $eCell.Value2 = $e - 1
    }
}
